$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '26.927.82'
$ws.Cells.Item(2, 5).Value = '  +0.00%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.671.25'
$ws.Cells.Item(3, 5).Value = '  +1.25%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.04%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '214.71'
$ws.Cells.Item(5, 5).Value = '  -0.04%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  +1.18%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.02%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '20.24'
$ws.Cells.Item(10, 5).Value = '  +0.39%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +1.65%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.906.96'
$ws.Cells.Item(12, 5).Value = '  +1.23%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '1.666.84'
$ws.Cells.Item(13, 5).Value = '  +1.01%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  +0.19%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  +1.36%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  +0.72%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '26.926.80'
$ws.Cells.Item(17, 5).Value = '  -0.03%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  +4.00%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '233.21'
$ws.Cells.Item(19, 5).Value = '  -0.93%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +0.16%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  -0.11%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +0.26%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '9.14'
$ws.Cells.Item(23, 5).Value = '  -1.85%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -1.84%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '146.35'
$ws.Cells.Item(25, 5).Value = '  +0.71%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.20%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '15.93'
$ws.Cells.Item(27, 5).Value = '  +0.93%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'BinanceUSD'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '1.00'
$ws.Cells.Item(28, 5).Value = '  +0.03%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'Stellar'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '0.112'
$ws.Cells.Item(29, 5).Value = '  -1.64%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '0.0498'
$ws.Cells.Item(30, 5).Value = '  +0.29%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +0.10%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '3.32'
$ws.Cells.Item(32, 5).Value = '  +0.62%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '1.453.75'
$ws.Cells.Item(33, 5).Value = '  -6.13%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +1.64%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +1.72%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +0.01%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.900'
$ws.Cells.Item(38, 5).Value = '  +1.11%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +0.76%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +13.03%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -4.12%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -0.02%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  +2.59%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '66.25'
$ws.Cells.Item(44, 5).Value = '  +0.79%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '1.811.41'
$ws.Cells.Item(45, 5).Value = '  +1.10%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.781'
$ws.Cells.Item(46, 5).Value = '  +0.84%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '90.65'
$ws.Cells.Item(47, 5).Value = '  +0.66%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +1.14%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +2.55%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +0.46%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '7.61'
$ws.Cells.Item(51, 5).Value = '  +0.14%  '
